# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" sheet (same layout as the other quarterly
#    sheets) right before the "总计" summary sheet.
# 2) Insert a new top row in "总计" summarizing the 2022-Q1 quarter and
#    renumber the index column underneath it.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2021-Q4")
$totalSheetBefore = $wb.Worksheets.Item("总计")

# Duplicate the most recent quarterly sheet so the new sheet starts with
# the same column layout / header styling, then rename + re-point it.
# NOTE: inserting a sheet shifts sheet positions around, and this engine's
# worksheet handles are position-based, so any sheet reference captured
# before the Copy() (including $totalSheetBefore itself) can silently
# start pointing at a different sheet afterwards. Re-resolve every sheet
# we still need by name once the sheet collection has settled.
$template.Copy($totalSheetBefore)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Item("总计")

# The template only has 5 fund rows (rows 2-5); 2022-Q1 needs 6 (rows 2-6),
# so insert one extra row and carry the formatting of the row above down
# onto it before filling in values.
$newSheet.Rows.Item(6).Insert()
$newSheet.Range("A5:H5").Copy()
$newSheet.Range("A6:H6").PasteSpecial(-4122)

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'968029"
$newSheet.Range("C2").Value = "恒生指数基金M类人民币（对冲）份额"
$newSheet.Range("D2").Value = "'25.09"
$newSheet.Range("E2").Value = "'97.94"
$newSheet.Range("F2").Value = "'8.76"
$newSheet.Range("G2").Value = "'2.1979"
$newSheet.Range("H2").Value = 1

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'539003"
$newSheet.Range("C3").Value = "建信富时100指数（QDII）人民币A"
$newSheet.Range("D3").Value = "'0.71"
$newSheet.Range("E3").Value = "'92.86"
$newSheet.Range("F3").Value = "'6.78"
$newSheet.Range("G3").Value = "'0.0481"
$newSheet.Range("H3").Value = 3

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'008707"
$newSheet.Range("C4").Value = "建信富时100指数（QDII）美元现汇A"
$newSheet.Range("D4").Value = "'0.71"
$newSheet.Range("E4").Value = "'92.86"
$newSheet.Range("F4").Value = "'6.78"
$newSheet.Range("G4").Value = "'0.0481"
$newSheet.Range("H4").Value = 3

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'008706"
$newSheet.Range("C5").Value = "建信富时100指数（QDII）人民币C"
$newSheet.Range("D5").Value = "'0.20"
$newSheet.Range("E5").Value = "'92.86"
$newSheet.Range("F5").Value = "'6.78"
$newSheet.Range("G5").Value = "'0.0136"
$newSheet.Range("H5").Value = 3

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'008708"
$newSheet.Range("C6").Value = "建信富时100指数（QDII）美元现汇C"
$newSheet.Range("D6").Value = "'0.20"
$newSheet.Range("E6").Value = "'92.86"
$newSheet.Range("F6").Value = "'6.78"
$newSheet.Range("G6").Value = "'0.0136"
$newSheet.Range("H6").Value = 3

# Now prepend the 2022-Q1 summary row to 总计, pushing the existing rows
# down and renumbering the index column (A) to stay 0-based.
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 5
$totalSheet.Range("D2").Value = 2.32

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
